# Update "想去人数" (F column) counts across sheets to match the
# regenerated gh-pages data output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 1015
$ws1.Range("F7").Value  = 607
$ws1.Range("F12").Value = 3031
$ws1.Range("F13").Value = 502
$ws1.Range("F14").Value = 1679
$ws1.Range("F18").Value = 1414
$ws1.Range("F21").Value = 1152
$ws1.Range("F22").Value = 20
$ws1.Range("F25").Value = 3574
$ws1.Range("F27").Value = 564
$ws1.Range("F28").Value = 1578

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 37
$ws2.Range("F7").Value = 4

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 800

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 800
$ws4.Range("F8").Value  = 37
$ws4.Range("F12").Value = 4
$ws4.Range("F16").Value = 1015
$ws4.Range("F18").Value = 607
$ws4.Range("F23").Value = 3031
$ws4.Range("F24").Value = 502
$ws4.Range("F25").Value = 1679
$ws4.Range("F29").Value = 1414
$ws4.Range("F34").Value = 1152
$ws4.Range("F35").Value = 20
$ws4.Range("F38").Value = 3574
$ws4.Range("F40").Value = 564
$ws4.Range("F41").Value = 1578
